$d = $word.ActiveDocument

# Position a collapsed range at the very start of the document body (before
# the existing first paragraph, which starts with "Brendan...").
$r = $d.Paragraphs.First.Range
$r.Collapse(1)

# Insert two new paragraphs ahead of the existing content:
#   1. A bold+italic warning paragraph.
#   2. An empty paragraph carrying an italic paragraph-mark.
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr><w:t>DON' + [char]0x2019 + 'T EDIT GO TO THE ONE IN THE ONEDRIVE</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:i/><w:iCs/></w:rPr></w:pPr></w:p>'

[void]$r.InsertXML($xml)
